$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.253.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.180.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.32'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.21'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -7.39%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.571'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.09'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '36.20'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -12.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0931'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.88'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.501.40'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.33'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.854'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.185.30'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.119.09'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.58'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.21%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.19%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.49%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.83%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.89'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.02'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.24'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.75'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0737'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.96'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.44'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.94%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.28'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +8.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.48'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -8.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.16'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -10.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.35'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.54'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.100'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.49%  '
